$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly data row inserted at the top of the data (row 9); the rest of
# the previously-existing rows (old 9..18) shift down by one (new 10..19).
# Column layout: D=Fecha, J=Volumen, K=Precio minimo, L=Precio maximo,
# M=Precio promedio ponderado, P=Precio $/Kg

$data = @(
    @{ Row = 9;  D = 44707; J = 30; K = 26000; L = 26000; M = 26000; P = 1733 }
    @{ Row = 10; D = 44428; J = 16; K = 25000; L = 26000; M = 25500; P = 1700 }
    @{ Row = 11; D = 44406; J = 25; K = 24000; L = 25000; M = 24520; P = 1635 }
    @{ Row = 12; D = 44680; J = 36; K = 24000; L = 25000; M = 24500; P = 1633 }
    @{ Row = 13; D = 44455; J = 18; K = 24000; L = 25000; M = 24500; P = 1633 }
    @{ Row = 14; D = 44385; J = 25; K = 14000; L = 15000; M = 14480; P = 965  }
    @{ Row = 15; D = 44421; J = 18; K = 24000; L = 25000; M = 24500; P = 1633 }
    @{ Row = 16; D = 44432; J = 34; K = 24000; L = 25000; M = 24500; P = 1633 }
    @{ Row = 17; D = 44329; J = 25; K = 23000; L = 23000; M = 23000; P = 1533 }
    @{ Row = 18; D = 44446; J = 34; K = 24000; L = 25000; M = 24500; P = 1633 }
    @{ Row = 19; D = 44705; J = 35; K = 26000; L = 26000; M = 26000; P = 1733 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 4).Value  = $entry.D   # D - Fecha
    $ws.Cells.Item($r, 10).Value = $entry.J   # J - Volumen
    $ws.Cells.Item($r, 11).Value = $entry.K   # K - Precio minimo
    $ws.Cells.Item($r, 12).Value = $entry.L   # L - Precio maximo
    $ws.Cells.Item($r, 13).Value = $entry.M   # M - Precio promedio ponderado
    $ws.Cells.Item($r, 16).Value = $entry.P   # P - Precio $/Kg
}
